$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '23.200.92'
Set-TextValue 'E2' '  -3.14%  '

Set-TextValue 'D3' '1.600.20'
Set-TextValue 'E3' '  -3.48%  '

Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  +0.50%  '

Set-TextValue 'E5' '  +0.25%  '

Set-TextValue 'D6' '301.89'
Set-TextValue 'E6' '  -2.46%  '

Set-TextValue 'D7' '0.3769'
Set-TextValue 'E7' '  -3.16%  '

Set-TextValue 'D8' '0.3674'
Set-TextValue 'E8' '  -4.25%  '

Set-TextValue 'D9' '49.05'
Set-TextValue 'E9' '  -3.93%  '

Set-TextValue 'D10' '1.005'
Set-TextValue 'E10' '  +0.54%  '

Set-TextValue 'D11' '1.276'
Set-TextValue 'E11' '  -5.83%  '

Set-TextValue 'D12' '0.08109'
Set-TextValue 'E12' '  -4.27%  '

Set-TextValue 'D13' '22.88'
Set-TextValue 'E13' '  -4.43%  '

Set-TextValue 'D14' '6.649'
Set-TextValue 'E14' '  -7.04%  '

Set-TextValue 'D15' '7.581'
Set-TextValue 'E15' '  -3.88%  '

Set-TextValue 'D16' '0.00001266'
Set-TextValue 'E16' '  -3.19%  '

Set-TextValue 'D17' '1.595.01'
Set-TextValue 'E17' '  -3.38%  '

Set-TextValue 'D18' '91.70'
Set-TextValue 'E18' '  -3.28%  '

Set-TextValue 'D19' '0.06829'
Set-TextValue 'E19' '  -2.50%  '

Set-TextValue 'D20' '18.50'
Set-TextValue 'E20' '  -6.55%  '

Set-TextValue 'D21' '6.616'
Set-TextValue 'E21' '  -4.27%  '

Set-TextValue 'D22' '1.003'
Set-TextValue 'E22' '  +0.11%  '

Set-TextValue 'D23' '13.16'
Set-TextValue 'E23' '  -3.24%  '

Set-TextValue 'D24' '23.199.83'
Set-TextValue 'E24' '  -3.10%  '

Set-TextValue 'B25' 'Toncoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D25' '2.363'
Set-TextValue 'E25' '  -4.80%  '

Set-TextValue 'B26' 'LidoDAOToken'
Set-TextValue 'C26' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D26' '2.967'
Set-TextValue 'E26' '  -2.81%  '

Set-TextValue 'D27' '21.16'
Set-TextValue 'E27' '  -4.18%  '

Set-TextValue 'D28' '150.84'
Set-TextValue 'E28' '  -1.03%  '

Set-TextValue 'D29' '5.298'
Set-TextValue 'E29' '  -2.08%  '

Set-TextValue 'D30' '132.37'
Set-TextValue 'E30' '  -4.93%  '

Set-TextValue 'D31' '2.472'
Set-TextValue 'E31' '  -0.91%  '

Set-TextValue 'D32' '7.125'
Set-TextValue 'E32' '  -8.43%  '

Set-TextValue 'D33' '1.774.50'
Set-TextValue 'E33' '  -3.15%  '

Set-TextValue 'D34' '0.9711'
Set-TextValue 'E34' '  -5.66%  '

Set-TextValue 'E35' '  -3.77%  '

Set-TextValue 'D36' '0.02787'
Set-TextValue 'E36' '  -5.67%  '

Set-TextValue 'D37' '6.302'
Set-TextValue 'E37' '  -5.41%  '

Set-TextValue 'B38' 'Algorand'
Set-TextValue 'C38' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D38' '0.2547'
Set-TextValue 'E38' '  -5.16%  '

Set-TextValue 'B39' 'FraxShare'
Set-TextValue 'C39' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D39' '10.17'
Set-TextValue 'E39' '  -7.49%  '

Set-TextValue 'D40' '0.08878'
Set-TextValue 'E40' '  -2.64%  '

Set-TextValue 'D41' '1.387'
Set-TextValue 'E41' '  -2.09%  '

Set-TextValue 'D42' '0.7182'
Set-TextValue 'E42' '  -4.76%  '

Set-TextValue 'D43' '12.84'
Set-TextValue 'E43' '  -4.68%  '

Set-TextValue 'D44' '16.30'
Set-TextValue 'E44' '  +0.77%  '

Set-TextValue 'D45' '0.6636'
Set-TextValue 'E45' '  -4.47%  '

Set-TextValue 'D46' '2.317'
Set-TextValue 'E46' '  -5.98%  '

Set-TextValue 'E47' '  +0.13%  '

Set-TextValue 'D48' '3.973'
Set-TextValue 'E48' '  -2.55%  '

Set-TextValue 'D49' '0.07998'
Set-TextValue 'E49' '  -3.41%  '

Set-TextValue 'D50' '131.71'
Set-TextValue 'E50' '  -1.95%  '

Set-TextValue 'D51' '1.172'
Set-TextValue 'E51' '  -4.86%  '
